$wb = $excel.ActiveWorkbook

$wsClient = $wb.Worksheets.Item("Client Info")
$wsSpace  = $wb.Worksheets.Item("Space Info")

# Move the "Community Membership" text from B3 to E3, and replace B3 with "Private Office"
$wsSpace.Range("E3").Value = "Community Membership"
$wsSpace.Range("B3").Value = "Private Office"

# Update the general comment section text
$wsSpace.Range("B4").Value = "General Comment Section Hello"

# Update selection on each sheet, then make "Space Info" the active tab/sheet
$wsClient.Range("B2").Select()
$wsSpace.Range("B5").Select()
$wsSpace.Activate()
